# Automatische test-sync: 2025-06-26 23:49:50
# Adds the 14th test-mail entry to the "Logs" sheet (row 46) and bumps the
# matching "Productinformatie" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 46 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$nl = [char]10

$logs.Range("A46").Value = "Kun je contact opnemen met de klant?"
$logs.Range("B46").Value = "mailmind.test@zohomail.eu"
$logs.Range("C46").Value = "Testmail #14: Kun je contact opnemen met de klant?"
$logs.Range("D46").Value = "Productinformatie"
$logs.Range("E46").Value = "Beste klantenservice," + $nl + "Ik heb zojuist een testmail verstuurd (Testmail #14) om te controleren of onze klanten goed bereikbaar zijn via e-mail. Kun je bevestigen of deze testmail succesvol is ontvangen en of jullie contact hebben opgenomen met de klant?" + $nl + "Met vriendelijke groet," + $nl + "[Naam]"
$logs.Range("F46").Value = "2025-06-26 23:49:29"
$logs.Range("G46").Value = "Ja"
$logs.Range("H46").Value = "Nee"
$logs.Range("I46").Value = "Ja"

# Re-collapse the row back to the default (non-custom) height: typing a
# multi-line value into E46 causes an Excel-like autofit to kick in (giving
# the row an explicit ht/customHeight), but the source row has none, so we
# explicitly autofit once more which clears the "custom" flag again.
$logs.Rows.Item(46).AutoFit()

# --- Logs sheet: extend conditional formatting ranges to cover row 46 ----
$logs.Range("D2:D45").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D46"))
$logs.Range("G2:G45").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G46"))
$logs.Range("H2:H45").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H46"))
$logs.Range("I2:I45").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I46"))

# --- Dashboard sheet: bump the Productinformatie count --------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B4").Value = 5
